# Update the localization status report:
#   - "Ready for handoff" -> "In Translation" for every Status cell
#   - Shrink the now-narrower "Status" columns to match the new text length
#     (Overview!E:F hold the per-language status, zh-cn!C and de-de!C hold
#     the "Status" table column)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2 and 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), rows 2 and 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C), rows 2 and 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Resize the columns that held the status text so they fit the shorter
#     string (mirrors what the reporting tool does when it regenerates the
#     sheet with narrower content) ---
$wsOverview.Columns(5).ColumnWidth = 12.5
$wsOverview.Columns(6).ColumnWidth = 12.5
$wsZhCn.Columns(3).ColumnWidth = 12.5
$wsDeDe.Columns(3).ColumnWidth = 12.5
